$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - bold/centered style like N1, values continue the sequence 14, 15
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Copy the style of N1 (header style) onto O1:P1
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-12, columns O and P
$ws.Range("O2").Value = -0.8459781767089031
$ws.Range("P2").Value = -0.6239004981304624

$ws.Range("O3").Value = -0.4600749389529301
$ws.Range("P3").Value = -0.3955961947093864

$ws.Range("O4").Value = 0.03294888504521559
$ws.Range("P4").Value = -0.01366074436541153

$ws.Range("O5").Value = 0.4243272129894631
$ws.Range("P5").Value = 0.386168490797377

$ws.Range("O6").Value = -0.3049156090296108
$ws.Range("P6").Value = -0.2967412013963857

$ws.Range("O7").Value = -0.1582994771796145
$ws.Range("P7").Value = -0.1583326034851112

$ws.Range("O8").Value = -0.4456834100645217
$ws.Range("P8").Value = -0.4384993929038604

$ws.Range("O9").Value = 0.003488422063327317
$ws.Range("P9").Value = 0.002889642554604474

$ws.Range("O10").Value = 0.00331207567487453
$ws.Range("P10").Value = 0.003984540936773856

$ws.Range("O11").Value = 0.009807123805076088
$ws.Range("P11").Value = 0.008748989017541843

$ws.Range("O12").Value = -0.01948723620516277
$ws.Range("P12").Value = -0.01916363945168371
